# Oracle to SQL Comparison.xlsx -- apply commit changes
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Oracle to SQL Samples")
$ws.Activate()

# ---------------------------------------------------------------------------
# Insert a new row at position 19 (everything at/after row 19 shifts down by 1)
# ---------------------------------------------------------------------------
$ws.Rows.Item(19).Insert()

# Copy formatting (styles) into the freshly inserted row 19 from rows that
# already carry the styles we need, so we match the target workbook exactly.
# A19/B19 -> style 11, C19 -> style 28 (same as row 21 A:C after the insert)
$ws.Range("A21:C21").Copy()
$ws.Range("A19:C19").PasteSpecial(-4122)   # xlPasteFormats
# D19 -> style 6 (same style as D34, which already used style 6 before insert shifted it)
$ws.Range("D34").Copy()
$ws.Range("D19").PasteSpecial(-4122)

$ws.Application.CutCopyMode = 0

# Row height for the new row (long multi-line text)
$ws.Rows.Item(19).RowHeight = 225.6

# New row 19 content
$ws.Range("A19").Value2 = "begin `ndrop index if exists [PROVIDER_LICENSE_IX0]  on [UA3_SCREENING].[PROVIDER_LICENSE];"
$ws.Range("B19").Value2 = "declare`n   V_COUNTER int;`nbegin`n  select count(*) into V_COUNTER`n  from sys.all_indexes`n  where OWNER = UA3_SCREENING`n    and TABLE_NAME = 'PROVIDER_LICENSE'`n    and INDEX_NAME = 'PROVIDER_LICENSE_IX0';`n  if V_COUNTER > 0 then`n    execute immediate 'drop index UA3_SCREENING.PROVIDER_LICENSE_IX0';`n  end if;  "
$ws.Range("C19").Value2 = "End with a GO statement"
$ws.Range("D19").Value2 = " "

# ---------------------------------------------------------------------------
# Row 20 (previously row 19, now shifted down by the insert above):
#   B20 gains the "execute immediate '...' " wrapper and its style becomes 11
# ---------------------------------------------------------------------------
$ws.Range("B3").Copy()
$ws.Range("B20").PasteSpecial(-4122)   # picks up style 11 for B20 (single cell -> single cell)
$ws.Application.CutCopyMode = 0
$ws.Range("B20").Value2 = "execute immediate 'create index PROVIDER_LICENSE_IX0 `n  on UA3_SCREENING.PROVIDER_LICENSE"

# ---------------------------------------------------------------------------
# Row 23 (previously row 22): B23 text ");" -> ")';"
# ---------------------------------------------------------------------------
$ws.Range("B23").Value2 = ")';"

# ---------------------------------------------------------------------------
# Row 24 (previously row 23, was blank): gains "end;" in both A24 and B24
# ---------------------------------------------------------------------------
$ws.Range("A24").Value2 = "end;"
$ws.Range("B24").Value2 = "end;"

# ---------------------------------------------------------------------------
# Row 25 (previously row 24): B25 gains "/"
# ---------------------------------------------------------------------------
$ws.Range("B25").Value2 = "/"

# ---------------------------------------------------------------------------
# Fix up the hyperlink that used to sit on (old) A60 -- now A61 after the
# row insertion -- the engine does not auto-shift the stored hyperlink ref.
# ---------------------------------------------------------------------------
$ws.Range("A60").Hyperlinks.Delete()
$ws.Hyperlinks.Add($ws.Range("A61"), "https://www.mssqltips.com/sqlservertutorial/2514/sql-server-insert-command/", "", "", "https://www.mssqltips.com/sqlservertutorial/2514/sql-server-insert-command/")

# ---------------------------------------------------------------------------
# Update the view: scrolled/selected area moved from B50 to B27
# ---------------------------------------------------------------------------
$ws.Range("B27").Select()
